# Update template analysis for plate reader growth curves
#
# Changes applied:
#  1. Rename the placeholder shared string used throughout the "strain"
#     sheet from "op_rep_selection" to "prom_T_marker".
#  2. Move the selection on the "strain" sheet from G6 to A10.
#  3. Give the "strain" sheet explicit column widths (A:F and H:L ~14.33,
#     G ~15.5).
#  4. Register an 8pt Calibri font (used for phonetic-guide metadata on
#     the "strain" sheet) without changing any cell's visible style.
#  5. Remove the now-unused "neg_selection" sheet. Excel automatically
#     re-targets the active tab to the last remaining sheet
#     ("pos_selection") once the trailing sheet is deleted.

$wb = $excel.ActiveWorkbook

# --- 1. strain sheet: relabel placeholder text -----------------------
$wsStrain = $wb.Worksheets.Item("strain")
$wsStrain.Range("A1:L8").Value = "prom_T_marker"

# --- 2. strain sheet: move selection to A10 ---------------------------
$wsStrain.Range("A10").Select()

# --- 3. strain sheet: explicit column widths --------------------------
$wsStrain.Columns("A:F").ColumnWidth = 13.5
$wsStrain.Columns("G:G").ColumnWidth = 14.7
$wsStrain.Columns("H:L").ColumnWidth = 13.5

# --- 4. register the 8pt Calibri font used by the phonetic guide ------
# Briefly apply it to a cell so it lands in the workbook's font table,
# then restore that cell's original 12pt size so no cell ends up with a
# non-default style.
$wsStrain.Range("A1").Font.Size = 8
$wsStrain.Range("A1").Font.Size = 12

# --- 5. drop the neg_selection sheet -----------------------------------
$excel.DisplayAlerts = $false
$wsNeg = $wb.Worksheets.Item("neg_selection")
$wsNeg.Delete()
$excel.DisplayAlerts = $true

Write-Output "edit complete"
